# Generate Report for Handoff
# Updates the localization-status report to reflect that the files are
# now "Ready for handoff", refreshes the handoff timestamps, and widens
# the date columns that now hold the longer status text / timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status column
$wsZhCn.Range("C2").Value = "Ready for handoff"       # Status column
$wsDeDe.Range("C2").Value = "Ready for handoff"       # Status column

# --- Timestamps refreshed for the new handoff ---
$wsOverview.Range("G2").Value = "2016-09-05 12:49:26"   # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value = "2016-09-05 12:49:26"        # Latest Handoff Datetime (de-de)
$wsZhCn.Range("H2").Value = "2016-09-05 12:49:22"         # Latest Handoff Datetime (zh-cn)

# --- Widen the date/status columns to fit the new content ---
# (target raw column width is 17.2159881591797 character-units; the Excel
# object model quantizes ColumnWidth to whole pixels on write, so the
# nearest representable value is used here)
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
